$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.989.81'
$ws.Range('E2').Value = '  +3.38%  '
$ws.Range('D3').Value = '3.451.46'
$ws.Range('E3').Value = '  +3.04%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'583.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.37%  '
$ws.Range('D6').Value = "'186.35"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.24%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '3.451.70'
$ws.Range('E8').Value = '  +3.29%  '
$ws.Range('D9').Value = "'0.999"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').Value = "'56.02"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.64%  '
$ws.Range('E13').Value = '  -1.50%  '
$ws.Range('D14').Value = "'9.37"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('D15').Value = '3.999.36'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = "'18.68"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').Value = '3.449.26'
$ws.Range('E17').Value = '  +3.26%  '
$ws.Range('D18').Value = '66.931.49'
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('D19').Value = "'12.09"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('D20').Value = "'0.118"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('E21').Value = '  +3.93%  '
$ws.Range('D22').Value = "'487.38"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.82%  '
$ws.Range('D23').Value = "'5.34"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.40%  '
$ws.Range('D24').Value = "'16.84"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +22.99%  '
$ws.Range('D25').Value = "'4.42"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.77%  '
$ws.Range('E26').Value = '  +3.15%  '
$ws.Range('E27').Value = '  +2.45%  '
$ws.Range('D28').Value = "'10.94"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.32%  '
$ws.Range('E29').Value = '  +4.73%  '
$ws.Range('D30').Value = "'31.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('D31').Value = "'7.17"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.78%  '
$ws.Range('D32').Value = "'598.47"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.20%  '
$ws.Range('D33').Value = "'11.72"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('E35').Value = '  +3.89%  '
$ws.Range('E36').Value = '  +7.21%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = "'36.61"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.04%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = "'0.384"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.84%  '
$ws.Range('D41').Value = '3.256.68'
$ws.Range('E41').Value = '  +5.45%  '
$ws.Range('D42').Value = '0.0₃0752'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('E43').Value = '  +5.83%  '
$ws.Range('D44').Value = "'0.0427"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('D45').Value = "'2.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +23.41%  '
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('D47').Value = "'3.23"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.55%  '
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +13.77%  '
$ws.Range('D50').Value = "'1.00"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = "'8.71"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.48%  '
